# Add a new "git ignore system (2)" entry row to the Git cheat-sheet (Sheet1),
# mirroring the existing "Git Ignore" row, and move the selection onto it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Populate the previously-empty row 28 with the new gitignore tip.
$ws.Range("A28").Value = "Git Ignore"
$ws.Range("B28").Value = "git ignore system (2)"
$ws.Range("C28").Value = "If some files are already under tracking, use below command to remove them from the working area:`n> git rm --cached -r target        //recursively remove files under target folder"

# Resize the row so the two-line note in C28 is fully visible.
$ws.Rows.Item(28).RowHeight = 26.25

# Move the active selection from C27 to the newly filled C28.
$ws.Range("C28").Select()
